$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.658.73"
$ws.Cells.Item(2, 5).Value = "  -7.18%  "
$ws.Cells.Item(3, 4).Value = "1.688.51"
$ws.Cells.Item(3, 5).Value = "  -6.45%  "
$ws.Cells.Item(4, 4).Value = "'1.006"
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$ws.Cells.Item(5, 4).Value = "'216.71"
$ws.Cells.Item(5, 5).Value = "  -6.40%  "
$ws.Cells.Item(6, 4).Value = "'1.007"
$ws.Cells.Item(6, 5).Value = "  +0.07%  "
$ws.Cells.Item(7, 4).Value = "'0.4948"
$ws.Cells.Item(7, 5).Value = "  -16.83%  "
$ws.Cells.Item(8, 4).Value = "'0.2584"
$ws.Cells.Item(8, 5).Value = "  -7.13%  "
$ws.Cells.Item(9, 4).Value = "'21.62"
$ws.Cells.Item(9, 5).Value = "  -7.73%  "
$ws.Cells.Item(10, 4).Value = "'0.06071"
$ws.Cells.Item(10, 5).Value = "  -11.36%  "
$ws.Cells.Item(11, 4).Value = "'0.07288"
$ws.Cells.Item(11, 5).Value = "  -3.44%  "
$ws.Cells.Item(12, 4).Value = "1.700.92"
$ws.Cells.Item(12, 5).Value = "  -5.85%  "
$ws.Cells.Item(13, 4).Value = "'4.426"
$ws.Cells.Item(13, 5).Value = "  -6.44%  "
$ws.Cells.Item(14, 4).Value = "1.918.81"
$ws.Cells.Item(14, 5).Value = "  -6.40%  "
$ws.Cells.Item(15, 4).Value = "'0.5711"
$ws.Cells.Item(15, 5).Value = "  -9.04%  "
$ws.Cells.Item(16, 4).Value = "'0.000008150"
$ws.Cells.Item(16, 5).Value = "  -12.09%  "
$ws.Cells.Item(17, 4).Value = "'64.56"
$ws.Cells.Item(17, 5).Value = "  -14.42%  "
$ws.Cells.Item(18, 4).Value = "26.677.40"
$ws.Cells.Item(18, 5).Value = "  -7.00%  "
$ws.Cells.Item(19, 4).Value = "'4.998"
$ws.Cells.Item(19, 5).Value = "  -8.79%  "
$ws.Cells.Item(20, 4).Value = "'1.006"
$ws.Cells.Item(20, 5).Value = "  +0.07%  "
$ws.Cells.Item(21, 5).Value = "  -6.55%  "
$ws.Cells.Item(22, 4).Value = "'181.73"
$ws.Cells.Item(22, 5).Value = "  -14.02%  "
$ws.Cells.Item(23, 4).Value = "'6.158"
$ws.Cells.Item(23, 5).Value = "  -10.31%  "
$ws.Cells.Item(24, 4).Value = "'1.007"
$ws.Cells.Item(24, 5).Value = "  +0.12%  "
$ws.Cells.Item(25, 4).Value = "'144.96"
$ws.Cells.Item(25, 5).Value = "  -6.15%  "
$ws.Cells.Item(26, 4).Value = "'7.521"
$ws.Cells.Item(26, 5).Value = "  -4.18%  "
$ws.Cells.Item(27, 4).Value = "'0.1126"
$ws.Cells.Item(27, 5).Value = "  -11.64%  "
$ws.Cells.Item(28, 4).Value = "'15.18"
$ws.Cells.Item(28, 5).Value = "  -7.59%  "
$ws.Cells.Item(29, 4).Value = "'1.312"
$ws.Cells.Item(29, 5).Value = "  -9.74%  "
$ws.Cells.Item(30, 4).Value = "'0.05566"
$ws.Cells.Item(30, 5).Value = "  -10.36%  "
$ws.Cells.Item(31, 4).Value = "'1.321"
$ws.Cells.Item(31, 5).Value = "  -7.00%  "
$ws.Cells.Item(32, 5).Value = "  -8.57%  "
$ws.Cells.Item(33, 4).Value = "'3.442"
$ws.Cells.Item(33, 5).Value = "  -8.32%  "
$ws.Cells.Item(34, 4).Value = "'1.641"
$ws.Cells.Item(34, 5).Value = "  -4.77%  "
$ws.Cells.Item(35, 4).Value = "'1.005"
$ws.Cells.Item(35, 5).Value = "  -4.90%  "
$ws.Cells.Item(36, 4).Value = "'2.402"
$ws.Cells.Item(36, 5).Value = "  -4.14%  "
$ws.Cells.Item(37, 4).Value = "'0.5838"
$ws.Cells.Item(37, 5).Value = "  -9.19%  "
$ws.Cells.Item(38, 4).Value = "'2.613"
$ws.Cells.Item(38, 5).Value = "  -4.13%  "
$ws.Cells.Item(39, 4).Value = "'0.01578"
$ws.Cells.Item(39, 5).Value = "  -7.75%  "
$ws.Cells.Item(40, 4).Value = "1.063.96"
$ws.Cells.Item(40, 5).Value = "  -7.02%  "
$ws.Cells.Item(41, 4).Value = "'5.855"
$ws.Cells.Item(41, 5).Value = "  -8.73%  "
$ws.Cells.Item(42, 4).Value = "'0.8468"
$ws.Cells.Item(42, 5).Value = "  -2.20%  "
$ws.Cells.Item(43, 4).Value = "'1.004"
$ws.Cells.Item(43, 5).Value = "  -0.21%  "
$ws.Cells.Item(44, 4).Value = "'97.85"
$ws.Cells.Item(44, 5).Value = "  -2.93%  "
$ws.Cells.Item(45, 4).Value = "1.847.68"
$ws.Cells.Item(45, 5).Value = "  -5.83%  "
$ws.Cells.Item(46, 4).Value = "'56.11"
$ws.Cells.Item(46, 5).Value = "  -7.44%  "
$ws.Cells.Item(47, 5).Value = "  -5.10%  "
$ws.Cells.Item(48, 4).Value = "'1.004"
$ws.Cells.Item(48, 5).Value = "  -0.35%  "
$ws.Cells.Item(49, 4).Value = "'8.045"
$ws.Cells.Item(49, 5).Value = "  -3.59%  "
$ws.Cells.Item(50, 5).Value = "  -3.74%  "
$ws.Cells.Item(51, 4).Value = "'0.05188"
$ws.Cells.Item(51, 5).Value = "  -5.07%  "
